$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date text strings from DD/MM/YYYY to DD-MM-YYYY for rows 3-21.
# Some of these (day <= 12) would otherwise be auto-parsed by Excel as
# real dates, so force them to stay literal text via the quote-prefix
# trick, then strip the resulting cell style back to Normal so the cell
# formatting is unchanged from before.
$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "'" + $dates[$row]
    $cell.Style = "Normal"
}

# Update attendance counts for specific rows: D and E flip 0 -> 1, H flips 1 -> 0
$rowsToFlip = @(11, 12, 14, 16)
foreach ($row in $rowsToFlip) {
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = 1
    $ws.Cells.Item($row, 8).Value = 0
}
